$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the two new sheets in the correct order: o_10, o_20, o_20_jumbled ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

# --- Header row (row 1): add 5th column 'evaluator_partial_correctness' to each sheet ---
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$ws1.Range("E1").Value = "evaluator_partial_correctness"

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)

# --- Row 2 data ---

# Sheet o_10: update prompt matrix, solution path, llm_response; add new evaluator_partial_correctness
$promptO10 = @'
 Given is the adjacency matrix for a weighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node I?
   A B C D E F G H I
 A 0 1 0 3 0 0 0 0 0
 B 1 0 2 0 2 0 0 0 0
 C 0 2 0 0 0 2 0 0 0
 D 3 0 0 0 1 0 2 0 0
 E 0 2 0 1 0 3 0 1 0
 F 0 0 2 0 3 0 0 0 1
 G 0 0 0 2 0 0 0 2 0
 H 0 0 0 0 1 0 2 0 1
 I 0 0 0 0 0 1 0 1 0

Solution: A -> B -> E -> H -> I
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 1 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 2 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 5 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 4 0 0 0 0 0 0 0 0
 E 3 0 0 0 0 2 0 0 3 0 0 0 0 0 0 0
 F 0 2 0 0 2 0 1 0 0 2 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 5 0 0 4 0 0 0 0 0
 H 0 0 0 4 0 0 5 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 3 0 0 0 0 0 0 0 3 0 0 0
 J 0 0 0 0 0 2 0 0 0 0 0 0 0 1 0 0
 K 0 0 0 0 0 0 4 0 0 0 0 0 0 0 2 0
 L 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 4
 M 0 0 0 0 0 0 0 0 3 0 0 0 0 4 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 4 0 4 0
 O 0 0 0 0 0 0 0 0 0 0 2 0 0 4 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 4 0 0 1 0
    
'@
$ws1.Range("A2").Value = $promptO10
$solO10 = @'
A -> B -> F -> G -> K -> O -> P
'@
$ws1.Range("B2").Value = $solO10
$llmO10 = @'
To find the least cost path from node A to node P, we can use Dijkstra's algorithm.
1. Initialize an empty set of visited nodes and a dictionary of distances. Set the distance of all nodes except A to infinity, and set the distance of A to 0.
2. While the set of visited nodes does not include P:
   a. Select the node with the smallest distance that has not been visited yet. Let's call this node "current".
   b. Mark current as visited.
   c. Update the distances of the neighboring nodes of current: for each neighbor, calculate the distance as the sum of the current distance and the cost of traveling from current to the neighbor. If this distance is smaller than the current distance of the neighbor, update the distance.
3. Once P has been visited, we know the least cost path from A to P is the path with the smallest distance. To find this path, we can backtrack from P to A using the updated distances and adjacency matrix:
   a. Initialize an empty path and set the current node to P.
   b. While the current node is not A:
      - Add the current node to the beginning of the path.
      - Find the neighbor of the current node that has the smallest distance.
      - Set the current node to this neighbor.
   c. Add A to the beginning of the path.
The final path we obtain will be the least cost path from A to P.
Applying this algorithm to the given adjacency matrix, we can find the least cost path from node A to node P.
'@
$ws1.Range("C2").Value = $llmO10
$ws1.Range("D2").Value = "Wrong"
$evalO10 = @'
Output: Not enough information provided
'@
$ws1.Range("E2").Value = $evalO10

# Sheet o_20
$promptO20 = @'
 Given is the adjacency matrix for a weighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 4 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 4 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 1 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 2
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 2 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 2 0 2 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0
Solution: A -> B -> C -> D -> I -> J -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 2 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 3 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 4 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 4 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 2 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 5 0 0 0 0 0 1 0 0 0 2 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 1 0 5 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 5 0 3 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 5 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 3 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 4 0 0 0 0 0 5
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 4 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 4 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 2 0 2 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 4
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0
    
'@
$ws2.Range("A2").Value = $promptO20
$solO20 = @'
A -> F -> K -> P -> Q -> V -> W -> X -> Y
'@
$ws2.Range("B2").Value = $solO20
$llmO20 = @'
The least cost path from node A to node Y is A -> F -> J -> O -> T -> Y.
'@
$ws2.Range("C2").Value = $llmO20
$ws2.Range("D2").Value = "Wrong"
$evalO20 = @'
Output: 2/6
'@
$ws2.Range("E2").Value = $evalO20

# Sheet o_20_jumbled
$promptO20j = @'
 Given is the adjacency matrix for a weighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 4 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 4 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 1 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 2
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 2 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 2 0 2 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0
Solution: A -> B -> C -> D -> I -> J -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 3 0 5 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 5 0 0 0 0 0 1 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 5 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 4 0 0 0 3 0 4 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 4 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 4 0 0 0 5 0 1 0 0 0 5 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 2 0 0 0 5 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 5 0 0 0 2 0 4 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 4 0 0 0 4 0 0 0 0 0 5 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0 0 5 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 2 0 0 0 0 0 5 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 4 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 4
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 4 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 2 0 5 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 2 0
    
'@
$ws3.Range("A2").Value = $promptO20j
$ws3.Range("B2").Value = $solO20
$llmO20j = @'
The least cost path from node A to node Y is A -> B -> C -> D -> I -> J -> O -> T -> Y.
'@
$ws3.Range("C2").Value = $llmO20j
$ws3.Range("D2").Value = "Wrong"
$evalO20j = @'
Output: 1/9
'@
$ws3.Range("E2").Value = $evalO20j

# Reset auto-grown row heights back to default (no wrap text in the source)
$ws1.Rows.Item(2).EntireRow.AutoFit()
$ws2.Rows.Item(2).EntireRow.AutoFit()
$ws3.Rows.Item(2).EntireRow.AutoFit()

# Re-activate o_10 so it remains the selected/active sheet
$ws1.Activate()
